$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture the current (pre-update) values of rows 451 and 452 ---
# These values will become the content of two brand-new rows that get
# inserted right after row 452 (the weekly update keeps the old
# observations and adds new ones on top).
$cols = 1..18
$row451 = @{}
$row452 = @{}
foreach ($c in $cols) {
    $row451[$c] = $ws.Cells.Item(451, $c).Value2
    $row452[$c] = $ws.Cells.Item(452, $c).Value2
}

# --- Step 2: insert two new blank rows before the old row 453 ---
# This pushes the old rows 453..467 down to 455..469, unchanged.
$ws.Rows.Item(453).Insert()
$ws.Rows.Item(453).Insert()

# --- Step 3: populate the two newly inserted rows (453, 454) with the
# values that used to be in rows 451 and 452 ---
foreach ($c in $cols) {
    $ws.Cells.Item(453, $c).Value = $row451[$c]
    $ws.Cells.Item(454, $c).Value = $row452[$c]
}

# --- Step 4: update rows 451 and 452 with their new values ---
$ws.Cells.Item(451, 4).Value = 44747
$ws.Cells.Item(451, 10).Value = 700
$ws.Cells.Item(451, 11).Value = 2000
$ws.Cells.Item(451, 12).Value = 2000
$ws.Cells.Item(451, 13).Value = 2000
$ws.Cells.Item(451, 16).Value = 2000

$ws.Cells.Item(452, 4).Value = 44747
$ws.Cells.Item(452, 10).Value = 700
$ws.Cells.Item(452, 11).Value = 2000
$ws.Cells.Item(452, 12).Value = 2000
$ws.Cells.Item(452, 13).Value = 2000
$ws.Cells.Item(452, 16).Value = 2000

Write-Host ("Used range rows now: " + $ws.UsedRange.Rows.Count)
